$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 62
$ws.Cells.Item($newRow, 1).Value = "2025-04-29 08:21:39"
$ws.Cells.Item($newRow, 2).Value = 172
